$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 318 (pushes existing rows 318-362 down to 321-365)
$ws.Rows("318:320").Insert()

# Common (unchanged across this product block) column values
$mercadoId = 9
$mercado = "Vega Central Mapocho de Santiago"
$region = "Metropolitana"
$codreg = 13
$tipo = "Fruta"
$productoId = 100107
$producto = "Otros"
$categoriaId = 100107002
$categoria = "Chirimoya"
$variedad = "Cultivar IV Región"

# New row data: Fecha, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, Unidad, Origen, PrecioKg, KgUnidad
$newRows = @(
    @{ Row = 318; Fecha = 45209; Calidad = "Especial"; Volumen = 290; PMin = 27000; PMax = 27000; PProm = 27000; Unidad = '$/bandeja 10 kilos'; Origen = "Provincia de Limarí"; PKg = 2700; KgUnidad = 10 },
    @{ Row = 319; Fecha = 45209; Calidad = "Primera";  Volumen = 300; PMin = 25000; PMax = 25000; PProm = 25000; Unidad = '$/bandeja 10 kilos'; Origen = "Provincia de Limarí"; PKg = 2500; KgUnidad = 10 },
    @{ Row = 320; Fecha = 45209; Calidad = "Segunda";  Volumen = 280; PMin = 22000; PMax = 22000; PProm = 22000; Unidad = '$/bandeja 10 kilos'; Origen = "Provincia de Limarí"; PKg = 2200; KgUnidad = 10 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
